$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "53.06.13.2021"
$ws.Range("A3").Value = "53.06.13.2020"
$ws.Range("A4").Value = "53.06.13.2019"
$ws.Range("A5").Value = "53.06.13.2018"
$ws.Range("A6").Value = "53.06.13.2017"
$ws.Range("A7").Value = "53.06.13.2016"
$ws.Range("A8").Value = "53.06.13.2015"
$ws.Range("A9").Value = "53.06.13.2014"

$ws.Range("G9").Select()
